$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 308 (high & close revised) ---
$ws.Cells.Item(308, 4).Value = 3.6919
$ws.Cells.Item(308, 6).Value = 3.6137

# --- Append new rows 309-311 with the same layout/format as the existing data rows ---
# Copy the date-cell formatting (style used in column A for data rows) down into the
# new rows first, then fill in the values.
$ws.Range("A308").Copy()
$ws.Range("A309:A311").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 309
$ws.Cells.Item(309, 1).Value = 45047.33333333334
$ws.Cells.Item(309, 2).Value = "FX_IDC:USDILS"
$ws.Cells.Item(309, 3).Value = 3.624
$ws.Cells.Item(309, 4).Value = 3.74944
$ws.Cells.Item(309, 5).Value = 3.6065
$ws.Cells.Item(309, 6).Value = 3.7154
$ws.Cells.Item(309, 7).Value = 0

# Row 310
$ws.Cells.Item(310, 1).Value = 45078.33333333334
$ws.Cells.Item(310, 2).Value = "FX_IDC:USDILS"
$ws.Cells.Item(310, 3).Value = 3.7155
$ws.Cells.Item(310, 4).Value = 3.7794
$ws.Cells.Item(310, 5).Value = 3.5401
$ws.Cells.Item(310, 6).Value = 3.6932
$ws.Cells.Item(310, 7).Value = 0

# Row 311
$ws.Cells.Item(311, 1).Value = 45110.33333333334
$ws.Cells.Item(311, 2).Value = "FX_IDC:USDILS"
$ws.Cells.Item(311, 3).Value = 3.7059
$ws.Cells.Item(311, 4).Value = 3.72694
$ws.Cells.Item(311, 5).Value = 3.6858
$ws.Cells.Item(311, 6).Value = 3.6949
$ws.Cells.Item(311, 7).Value = 0
